# "Generate Report for Handoff"
#
# The handoff transform for the 41b545a2-...md source file failed this run, so the
# localization-status report needs to reflect that:
#   - Status changes from "Ready for handoff" to "Handoff transform failed"
#     (shown on the Overview summary sheet as well as on each per-locale sheet)
#   - the "Latest Handoff File" link/value is cleared out, since no handoff file was
#     produced this time
#   - "Latest Handoff Datetime", "Latest Handback DateTime" and "Handoff Reason" are
#     reset back to the same "nothing has happened yet" defaults used by the ignored
#     .localization-config row

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Handoff transform failed"
$wsOverview.Range("C2").Value = "Handoff transform failed"

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Drop the "Latest Handoff File" hyperlink + value in row 2 (column C)
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq '$C$2') {
            $hl.Delete()
        }
    }
    $ws.Range("C2").Clear()

    # Update the Status text for the failed handoff
    $ws.Range("B2").Value = "Handoff transform failed"

    # Reset the datetime / reason columns back to their "not yet happened" defaults,
    # matching row 3 (the ignored .localization-config entry)
    $ws.Range("D2").Value = $ws.Range("D3").Value()
    $ws.Range("G2").Value = $ws.Range("G3").Value()
    $ws.Range("H2").Value = $ws.Range("H3").Value()
}
